$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - read_csv_to_df test
$ws.Range("A4").Value = "s3"
$ws.Range("B4").Value = "read_csv_to_df"
$ws.Range("C4").Value = "test_read_csv_to_df_ok"
$ws.Range("D4").Value = "Test the read_csv_to_df method for reading 1 .csv file from the mocked s3 bucket"
$ws.Range("E4").Value = "Mock s3 bucket`nRead csv file from it into DataFrame"
$ws.Range("E4").WrapText = $true
$ws.Range("F4").Value = "test.csv"
$ws.Range("G4").Value = "DataFrame with the content of the csv file"
$ws.Range("G4").WrapText = $true

# Row 5 - write_df_to_s3 empty DataFrame test
$ws.Range("A5").Value = "s3"
$ws.Range("B5").Value = "write_df_to_s3"
$ws.Range("C5").Value = "test_write_df_to_s3_empty"
$ws.Range("D5").Value = "Test the write_df_to_s3 method within an empty DataFrame as input"
$ws.Range("F5").Value = "empty DataFrame"
$ws.Range("G5").Value = "None"

# Row 6 - write_df_to_s3 csv test
$ws.Range("A6").Value = "s3"
$ws.Range("B6").Value = "write_df_to_s3"
$ws.Range("C6").Value = "test_write_df_to_s3_csv"
$ws.Range("D6").Value = "Tests if writting a csv file to write_df_to_s3 method succeeds"
$ws.Range("F6").Value = "test.csv"

# Row 7 - write_df_to_s3 parquet test
$ws.Range("A7").Value = "s3"
$ws.Range("B7").Value = "write_df_to_s3"
$ws.Range("C7").Value = "test_write_df_to_s3_parquet"
$ws.Range("D7").Value = "Tests if writting a parquet file to write_df_to_s3 method succeeds"
$ws.Range("F7").Value = "test.parquet"

# Row heights to match the other data rows
$ws.Rows.Item(4).RowHeight = 43.2
$ws.Rows.Item(5).RowHeight = 43.2
$ws.Rows.Item(6).RowHeight = 43.2
$ws.Rows.Item(7).RowHeight = 43.2

# Widen column F to fit the new content
$ws.Columns.Item(6).ColumnWidth = 14.83

# Move the selection the way it ends up after the edits
[void]$ws.Range("F6").Select()
